$p = $ppt.ActivePresentation

# ---------------------------------------------------------------------------
# Insert two new "Software Requirements" slides immediately before the
# existing "Server Implementation" (threads) slide, which currently sits at
# index 15. We build them by duplicating the "Server Implementation"
# (cached-list) slide at index 16 -- it already uses the same layout as the
# surrounding deck and its placeholders carry plain (non-autofit) bodyPr and
# existing en-US runs we can safely overwrite -- then relocating the
# duplicates into place.
# ---------------------------------------------------------------------------
$cachedListSlide = $p.Slides.Item(16)
$reqSlide1 = $cachedListSlide.Duplicate()
$reqSlide1.MoveTo(15)

$cachedListSlide2 = $p.Slides.Item(17)
$reqSlide2 = $cachedListSlide2.Duplicate()
$reqSlide2.MoveTo(16)

# --- Slide 15: "Software Requirements" / Functional requirements ----------
$body1 = $reqSlide1.Shapes.Item(1).TextFrame.TextRange
$body1.Text = "Functional:`rServer should be robust and be resilient to failure`rServer log messages should be clear on activity`rDevice crashes should not corrupt any part of server`rDevices may only have one “owner” at any given instance`rClients can actively control only one device any given instance`rServer must be secure against unwarranted input"
$body1.Paragraphs(2,1).IndentLevel = 2
$body1.Paragraphs(3,1).IndentLevel = 2
$body1.Paragraphs(4,1).IndentLevel = 2
$body1.Paragraphs(5,1).IndentLevel = 2
$body1.Paragraphs(6,1).IndentLevel = 2
$body1.Paragraphs(7,1).IndentLevel = 2

$title1 = $reqSlide1.Shapes.Item(2).TextFrame.TextRange
$title1.Text = "Software Requirements"

# --- Slide 16: "Software Requirements" / Non-functional requirements ------
$body2 = $reqSlide2.Shapes.Item(1).TextFrame.TextRange
$body2.Text = "Non-functional:`rReal-time devices require near real-time feedback`rDevices should have minimal setup to boot up and connect to server`rSystem should be responsive under any amount of stress`rServer deployment should be straight forward`r`r"
$body2.Paragraphs(2,1).IndentLevel = 2
$body2.Paragraphs(3,1).IndentLevel = 2
$body2.Paragraphs(4,1).IndentLevel = 2
$body2.Paragraphs(5,1).IndentLevel = 2
$body2.Paragraphs(6,1).IndentLevel = 2
$body2.Paragraphs(7,1).IndentLevel = 2

$title2 = $reqSlide2.Shapes.Item(2).TextFrame.TextRange
$title2.Text = "Software Requirements"
